$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("API")

# The "Brand" block in columns E:F (rows 4-20) needs a new row inserted at
# the top (row 4) for "View single brand" / "localhost:3000/viewBrandById",
# pushing the existing View/Add/Edit/Delete Brand + Discount/DiscountType/Cart
# blocks down by one row. Columns B:C and H:I must stay untouched.
#
# Range.Insert on this host shifts whole rows, so instead we manually walk
# the E:F column band bottom-up and copy each row's value+format down one
# row; processing bottom-up avoids clobbering data we still need to read.
for ($r = 20; $r -ge 4; $r--) {
    $srcE = $ws.Cells.Item($r, 5)
    $srcF = $ws.Cells.Item($r, 6)
    $dstE = $ws.Cells.Item($r + 1, 5)
    $dstF = $ws.Cells.Item($r + 1, 6)

    $vE = $srcE.Value()
    $vF = $srcF.Value()

    $ws.Range($srcE, $srcF).Copy()
    $ws.Range($dstE, $dstF).PasteSpecial(-4122) # xlPasteFormats

    $dstE.Value = $vE
    $dstF.Value = $vF
}

# New row 4: "View single brand" styled like the other yellow header cells
# (E2:F2 / E3:F3) above it.
$ws.Range("E2:F2").Copy()
$ws.Range("E4:F4").PasteSpecial(-4122) # xlPasteFormats
$ws.Cells.Item(4, 5).Value = "View single brand"
$ws.Cells.Item(4, 6).Value = "localhost:3000/viewBrandById"

$excel.CutCopyMode = 0

# Update the API sheet's remembered selection.
$ws.Activate() | Out-Null
$ws.Range("F4").Select() | Out-Null

# Update the Table sheet's remembered selection (it stays the active tab).
$wsTable = $wb.Worksheets.Item("Table")
$wsTable.Activate() | Out-Null
$wsTable.Range("H16").Select() | Out-Null
